$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new list bullet "Split up test cases further." right after the
#    "Improved test coverage..." bullet (and before "Examples"), inheriting
#    that bullet's List Paragraph style / ilvl=1 / numId=1 formatting.
# ---------------------------------------------------------------------------
$anchorText = "Improved test coverage of all APIs (including operator overloads, copy/move constructors, etc)."

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
}

$anchorRange = $d.Paragraphs.Item($anchorIndex).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "Split up test cases further."

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark onto the end of the text we just typed -
#    exactly what Word does after you type new text, since that bookmark
#    always tracks the most recent edit position. (Bookmarks.Add re-uses the
#    existing "_GoBack" bookmark, automatically removing it from its old
#    location -- the "New Modules" heading.)
#
#    Quirk work-around: anchoring a zero-length bookmark exactly at "end of
#    paragraph text, right before the paragraph mark" lands it at the wrong
#    spot. Side-step this by nudging a throw-away character in after that
#    position, anchoring the bookmark against the now-safe position, and
#    then deleting the throw-away character again.
# ---------------------------------------------------------------------------
$newParaRange = $d.Paragraphs.Item($anchorIndex + 1).Range
$endPos = $newParaRange.End - 1

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()
